# Update the two-digit-division answer table to the newly generated set of
# problems/answers. Most cells are simple text substitutions; the last
# populated row also loses its first cell (old "18÷7=2, 4") and gains a new
# cell at the end (new "15÷2=7, 1"), so that row is rewritten cell-by-cell
# instead of via Find/Replace.

$d = $word.ActiveDocument

function Replace-Answer($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Row 1
Replace-Answer "52÷4=13, 0" "65÷7=9, 2"
Replace-Answer "44÷7=6, 2" "45÷4=11, 1"
Replace-Answer "83÷6=13, 5" "46÷8=5, 6"
Replace-Answer "80÷7=11, 3" "10÷8=1, 2"
Replace-Answer "51÷3=17, 0" "16÷4=4, 0"

# Row 2
Replace-Answer "14÷2=7, 0" "41÷9=4, 5"
Replace-Answer "16÷3=5, 1" "58÷8=7, 2"
Replace-Answer "17÷8=2, 1" "43÷2=21, 1"
Replace-Answer "43÷6=7, 1" "54÷8=6, 6"
Replace-Answer "12÷3=4, 0" "13÷4=3, 1"

# Row 3
Replace-Answer "47÷8=5, 7" "53÷7=7, 4"
Replace-Answer "15÷4=3, 3" "15÷7=2, 1"
Replace-Answer "58÷7=8, 2" "91÷7=13, 0"
Replace-Answer "96÷6=16, 0" "39÷2=19, 1"
Replace-Answer "89÷2=44, 1" "91÷8=11, 3"

# Row 4
Replace-Answer "74÷2=37, 0" "75÷9=8, 3"
Replace-Answer "71÷6=11, 5" "60÷3=20, 0"
Replace-Answer "80÷5=16, 0" "60÷9=6, 6"
Replace-Answer "83÷5=16, 3" "86÷5=17, 1"
Replace-Answer "62÷8=7, 6" "54÷8=6, 6"

# Row 5: the first cell ("18÷7=2, 4") is removed entirely, shifting the
# remaining four cells left by one and adding a brand-new fifth cell
# ("15÷2=7, 1"). Net effect on the five visible cell slots of that row:
$t = $d.Tables.Item(1)
$row = 17
$t.Cell($row, 1).Range.Text = "14÷4=3, 2"
$t.Cell($row, 2).Range.Text = "88÷9=9, 7"
$t.Cell($row, 3).Range.Text = "90÷6=15, 0"
$t.Cell($row, 4).Range.Text = "73÷5=14, 3"
$t.Cell($row, 5).Range.Text = "15÷2=7, 1"
